$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4000, 4000),
    @(4000, 1000),
    @(5000, 4000),
    @(500, 500),
    @(500, 500),
    @(500, 500),
    @(500, 500)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
